# update scripts with new tpm
# Refresh the Mdk -> Tspan1 ligand-receptor table with recalculated TPM-based
# values, extending the sender/target cluster grid to also include the
# "Resolving-Mac" cluster as a possible target (rows 2-17, columns A-T).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then the 16 numeric metric columns (E through T).
$data = @(
    ("ECs", "Mdk", "Tspan1", "ECs", 3, 1, 2.180165333333334, 6.540496, 0.01970539991828544, 0.01970539991828544, 1, 0.3333333333333333, 0.1992763333333334, 0.5978290000000001, 0.08623076327255404, 0.08623076327255402, 0.4344553536871112, 3.910098183184, 0.001699211675544677, 0.001699211675544677),
    ("ECs", "Mdk", "Tspan1", "FAPs", 3, 1, 2.180165333333334, 6.540496, 0.01970539991828544, 0.01970539991828544, 3, 1, 1.489416666666667, 4.46825, 0.6444996947163647, 0.6444996947163646, 3.247174583555556, 29.224571252, 0.01270012423159884, 0.01270012423159884),
    ("ECs", "Mdk", "Tspan1", "MuSCs", 3, 1, 2.180165333333334, 6.540496, 0.01970539991828544, 0.01970539991828544, 2, 0.6666666666666666, 0.572399, 1.717197, 0.2476882319180568, 0.2476882319180568, 1.247924456634667, 11.231320109712, 0.00488079566499834, 0.00488079566499834),
    ("ECs", "Mdk", "Tspan1", "Resolving-Mac", 3, 1, 2.180165333333334, 6.540496, 0.01970539991828544, 0.01970539991828544, 1, 0.3333333333333333, 0.04987366666666667, 0.149621, 0.0215813100930246, 0.0215813100930246, 0.1087328391128889, 0.978595552016, 0.0004252683461435797, 0.0004252683461435797),
    ("FAPs", "Mdk", "Tspan1", "ECs", 3, 1, 81.17653533333333, 243.529606, 0.733713204346044, 0.7337132043460441, 1, 0.3333333333333333, 0.1992763333333334, 0.5978290000000001, 0.08623076327255404, 0.08623076327255402, 16.17656231393045, 145.589060825374, 0.0632686496339108, 0.0632686496339108),
    ("FAPs", "Mdk", "Tspan1", "FAPs", 3, 1, 81.17653533333333, 243.529606, 0.733713204346044, 0.7337132043460441, 3, 1, 1.489416666666667, 4.46825, 0.6444996947163647, 0.6444996947163646, 120.9056846677222, 1088.1511620095, 0.4728779362103911, 0.4728779362103911),
    ("FAPs", "Mdk", "Tspan1", "MuSCs", 3, 1, 81.17653533333333, 243.529606, 0.733713204346044, 0.7337132043460441, 2, 0.6666666666666666, 0.572399, 1.717197, 0.2476882319180568, 0.2476882319180568, 46.46536764826467, 418.188308834382, 0.1817321263194035, 0.1817321263194036),
    ("FAPs", "Mdk", "Tspan1", "Resolving-Mac", 3, 1, 81.17653533333333, 243.529606, 0.733713204346044, 0.7337132043460441, 1, 0.3333333333333333, 0.04987366666666667, 0.149621, 0.0215813100930246, 0.0215813100930246, 4.048571464369556, 36.437143179326, 0.01583449218233871, 0.01583449218233871),
    ("MuSCs", "Mdk", "Tspan1", "ECs", 3, 1, 25.672264, 77.016792, 0.2320384702908474, 0.2320384702908474, 1, 0.3333333333333333, 0.1992763333333334, 0.5978290000000001, 0.08623076327255404, 0.08623076327255402, 5.115874638285334, 46.042871744568, 0.02000885440177563, 0.02000885440177562),
    ("MuSCs", "Mdk", "Tspan1", "FAPs", 3, 1, 25.672264, 77.016792, 0.2320384702908474, 0.2320384702908474, 3, 1, 1.489416666666667, 4.46825, 0.6444996947163647, 0.6444996947163646, 38.23669787266667, 344.130280854, 0.1495487232649034, 0.1495487232649034),
    ("MuSCs", "Mdk", "Tspan1", "MuSCs", 3, 1, 25.672264, 77.016792, 0.2320384702908474, 0.2320384702908474, 2, 0.6666666666666666, 0.572399, 1.717197, 0.2476882319180568, 0.2476882319180568, 14.694778241336, 132.253004172024, 0.05747319844331053, 0.05747319844331053),
    ("MuSCs", "Mdk", "Tspan1", "Resolving-Mac", 3, 1, 25.672264, 77.016792, 0.2320384702908474, 0.2320384702908474, 1, 0.3333333333333333, 0.04987366666666667, 0.149621, 0.0215813100930246, 0.0215813100930246, 1.280369937314667, 11.523329435832, 0.005007694180857855, 0.005007694180857855),
    ("Resolving-Mac", "Mdk", "Tspan1", "ECs", 3, 1, 1.608999666666667, 4.826999, 0.01454292544482312, 0.01454292544482312, 1, 0.3333333333333333, 0.1992763333333334, 0.5978290000000001, 0.08623076327255404, 0.08623076327255402, 0.3206355539078889, 2.885719985171, 0.001254047561322945, 0.001254047561322945),
    ("Resolving-Mac", "Mdk", "Tspan1", "FAPs", 3, 1, 1.608999666666667, 4.826999, 0.01454292544482312, 0.01454292544482312, 3, 1, 1.489416666666667, 4.46825, 0.6444996947163647, 0.6444996947163646, 2.396470920194445, 21.56823828175, 0.009372911009471355, 0.009372911009471354),
    ("Resolving-Mac", "Mdk", "Tspan1", "MuSCs", 3, 1, 1.608999666666667, 4.826999, 0.01454292544482312, 0.01454292544482312, 2, 0.6666666666666666, 0.572399, 1.717197, 0.2476882319180568, 0.2476882319180568, 0.9209898002003333, 8.288908201803, 0.003602111490344359, 0.003602111490344359),
    ("Resolving-Mac", "Mdk", "Tspan1", "Resolving-Mac", 3, 1, 1.608999666666667, 4.826999, 0.01454292544482312, 0.01454292544482312, 1, 0.3333333333333333, 0.04987366666666667, 0.149621, 0.0215813100930246, 0.0215813100930246, 0.08024671304211112, 0.722220417379, 0.0003138553836844656, 0.0003138553836844656)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}
